$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = $ws.Range("D6:D7")
$r.Style = "Good"
$r.HorizontalAlignment = -4108
